$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price-range values in column B (order matters for shared-string layout)
$ws.Range("B9").Value = "1.139/1.161"
$ws.Range("B8").Value = "1.114/1.150"

# Clear the now-unused D column values
$ws.Range("D9").ClearContents()

$ws.Range("B18").Value = "5.005/5.029"
$ws.Range("D18").ClearContents()

# B12 keeps the same displayed text, but update it anyway to mirror the source edit
$ws.Range("B12").Value = "300ETF（510300）"

# Update the selected cell in the sheet view
$ws.Range("B22").Select()
